# Weekly update: insert two new price records at the top of the data
# block (rows 121-122), pushing all existing records down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the two new rows (everything currently on/after row 121
# shifts down by one for each Insert() call).
$ws.Rows.Item(121).Insert()
$ws.Rows.Item(122).Insert()

# New row 121
$ws.Cells.Item(121, 1).Value = 10
$ws.Cells.Item(121, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(121, 3).Value = "La Araucanía"
$ws.Cells.Item(121, 4).Value = [DateTime]"2022-07-20"
$ws.Cells.Item(121, 5).Value = 9
$ws.Cells.Item(121, 6).Value = "Fruta"
$ws.Cells.Item(121, 7).Value = 100104
$ws.Cells.Item(121, 8).Value = "Frutos de pepita"
$ws.Cells.Item(121, 9).Value = 100104003
$ws.Cells.Item(121, 10).Value = "Membrillo"
$ws.Cells.Item(121, 11).Value = "Champion"
$ws.Cells.Item(121, 12).Value = "Primera"
$ws.Cells.Item(121, 13).Value = 100
$ws.Cells.Item(121, 14).Value = 10000
$ws.Cells.Item(121, 15).Value = 10000
$ws.Cells.Item(121, 16).Value = 10000
$ws.Cells.Item(121, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(121, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(121, 19).Value = 556
$ws.Cells.Item(121, 20).Value = 18

# New row 122
$ws.Cells.Item(122, 1).Value = 10
$ws.Cells.Item(122, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(122, 3).Value = "La Araucanía"
$ws.Cells.Item(122, 4).Value = [DateTime]"2022-07-20"
$ws.Cells.Item(122, 5).Value = 9
$ws.Cells.Item(122, 6).Value = "Fruta"
$ws.Cells.Item(122, 7).Value = 100104
$ws.Cells.Item(122, 8).Value = "Frutos de pepita"
$ws.Cells.Item(122, 9).Value = 100104003
$ws.Cells.Item(122, 10).Value = "Membrillo"
$ws.Cells.Item(122, 11).Value = "Champion"
$ws.Cells.Item(122, 12).Value = "Primera"
$ws.Cells.Item(122, 13).Value = 3
$ws.Cells.Item(122, 14).Value = 200000
$ws.Cells.Item(122, 15).Value = 200000
$ws.Cells.Item(122, 16).Value = 200000
$ws.Cells.Item(122, 17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(122, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(122, 19).Value = 444
$ws.Cells.Item(122, 20).Value = 450
